$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "gen_3" worksheet right before "Función objetivo"
# ---------------------------------------------------------------------------
$funcSheet = $wb.Worksheets.Item("Función objetivo")
$gen3 = $wb.Worksheets.Add($funcSheet, $null)
$gen3.Name = "gen_3"

# ---------------------------------------------------------------------------
# 2. Update gen_0 (sheet1) rows 2-11
# ---------------------------------------------------------------------------
$gen0 = $wb.Worksheets.Item("gen_0")

$gen0.Range("A2").Value = -0.9682539682539683
$gen0.Range("B2").Value = 1
$gen0.Range("C2").Value = "000001"
$gen0.Range("D2").Value = 0.9375157470395565
$gen0.Range("F2").Value = "001100"
$gen0.Range("G2").Value = "['110100' '001000']"
$gen0.Range("I2").Value = -0.6190476190476191
$gen0.Range("J2").Value = 12
$gen0.Range("K2").Value = "001100"
$gen0.Range("L2").Value = 0.3832199546485261

$gen0.Range("A3").Value = -0.746031746031746
$gen0.Range("B3").Value = 8
$gen0.Range("C3").Value = "001000"
$gen0.Range("F3").Value = "110000"
$gen0.Range("G3").Value = "['110100' '001000']"
$gen0.Range("I3").Value = 0.5238095238095237
$gen0.Range("J3").Value = 48
$gen0.Range("K3").Value = "110000"
$gen0.Range("L3").Value = 0.27437641723356

$gen0.Range("A4").Value = -0.7142857142857143
$gen0.Range("B4").Value = 9
$gen0.Range("C4").Value = "001001"
$gen0.Range("D4").Value = 0.5102040816326531
$gen0.Range("F4").Value = "001110"
$gen0.Range("G4").Value = "['110110' '001000']"
$gen0.Range("I4").Value = -0.5555555555555556
$gen0.Range("J4").Value = 14
$gen0.Range("K4").Value = "001110"
$gen0.Range("L4").Value = 0.308641975308642

$gen0.Range("A5").Value = 0.7142857142857142
$gen0.Range("B5").Value = 54
$gen0.Range("C5").Value = "110110"
$gen0.Range("D5").Value = 0.510204081632653
$gen0.Range("F5").Value = "110000"
$gen0.Range("G5").Value = "['110110' '001000']"
$gen0.Range("I5").Value = 0.5238095238095237
$gen0.Range("J5").Value = 48
$gen0.Range("K5").Value = "110000"
$gen0.Range("L5").Value = 0.27437641723356

$gen0.Range("A6").Value = 0.6507936507936507
$gen0.Range("B6").Value = 52
$gen0.Range("C6").Value = "110100"
$gen0.Range("D6").Value = 0.4235323759133282
$gen0.Range("F6").Value = "110001"
$gen0.Range("G6").Value = "['000001' '110100']"
$gen0.Range("I6").Value = 0.5555555555555554
$gen0.Range("J6").Value = 49
$gen0.Range("K6").Value = "110001"
$gen0.Range("L6").Value = 0.3086419753086417

$gen0.Range("A7").Value = 0.5238095238095237
$gen0.Range("B7").Value = 48
$gen0.Range("C7").Value = "110000"
$gen0.Range("D7").Value = 0.27437641723356
$gen0.Range("F7").Value = "000100"
$gen0.Range("G7").Value = "['000001' '110100']"
$gen0.Range("I7").Value = -0.873015873015873
$gen0.Range("J7").Value = 4
$gen0.Range("K7").Value = "000100"
$gen0.Range("L7").Value = 0.762156714537667

$gen0.Range("A8").Value = -0.3650793650793651
$gen0.Range("B8").Value = 20
$gen0.Range("C8").Value = "010100"
$gen0.Range("D8").Value = 0.1332829428067523
$gen0.Range("F8").Value = "001000"
$gen0.Range("G8").Value = "['001000' '001001']"
$gen0.Range("I8").Value = -0.746031746031746
$gen0.Range("J8").Value = 8
$gen0.Range("K8").Value = "001000"
$gen0.Range("L8").Value = 0.5565633660871756

$gen0.Range("A9").Value = 0.2698412698412698
$gen0.Range("B9").Value = 40
$gen0.Range("C9").Value = "101000"
$gen0.Range("D9").Value = 0.07281431090954897
$gen0.Range("F9").Value = "001001"
$gen0.Range("G9").Value = "['001000' '001001']"
$gen0.Range("I9").Value = -0.7142857142857143
$gen0.Range("J9").Value = 9
$gen0.Range("K9").Value = "001001"
$gen0.Range("L9").Value = 0.5102040816326531

$gen0.Range("A10").Value = 0.2063492063492063
$gen0.Range("B10").Value = 38
$gen0.Range("C10").Value = "100110"
$gen0.Range("D10").Value = 0.04257999496094732
$gen0.Range("F10").Value = "001001"
$gen0.Range("G10").Value = "['001001' '001001']"
$gen0.Range("I10").Value = -0.7142857142857143
$gen0.Range("J10").Value = 9
$gen0.Range("K10").Value = "001001"
$gen0.Range("L10").Value = 0.5102040816326531

$gen0.Range("A11").Value = 0.1746031746031744
$gen0.Range("B11").Value = 37
$gen0.Range("C11").Value = "100101"
$gen0.Range("D11").Value = 0.03048626858150661
$gen0.Range("F11").Value = "001001"
$gen0.Range("G11").Value = "['001001' '001001']"
$gen0.Range("I11").Value = -0.7142857142857143
$gen0.Range("J11").Value = 9
$gen0.Range("K11").Value = "001001"
$gen0.Range("L11").Value = 0.5102040816326531

# ---------------------------------------------------------------------------
# 3. Update gen_1 (sheet2) rows 2-11
# ---------------------------------------------------------------------------
$gen1 = $wb.Worksheets.Item("gen_1")

$gen1.Range("A2").Value = -0.9682539682539683
$gen1.Range("B2").Value = 1
$gen1.Range("C2").Value = "000001"
$gen1.Range("D2").Value = 0.9375157470395565
$gen1.Range("F2").Value = "001000"
$gen1.Range("G2").Value = "['001000' '001000']"
$gen1.Range("I2").Value = -0.746031746031746
$gen1.Range("J2").Value = 8
$gen1.Range("K2").Value = "001000"
$gen1.Range("L2").Value = 0.5565633660871756

$gen1.Range("A3").Value = -0.873015873015873
$gen1.Range("B3").Value = 4
$gen1.Range("C3").Value = "000100"
$gen1.Range("D3").Value = 0.762156714537667
$gen1.Range("F3").Value = "001000"
$gen1.Range("G3").Value = "['001000' '001000']"
$gen1.Range("I3").Value = -0.746031746031746
$gen1.Range("J3").Value = 8
$gen1.Range("K3").Value = "001000"
$gen1.Range("L3").Value = 0.5565633660871756

$gen1.Range("A4").Value = -0.746031746031746
$gen1.Range("B4").Value = 8
$gen1.Range("C4").Value = "001000"
$gen1.Range("D4").Value = 0.5565633660871756
$gen1.Range("F4").Value = "110000"
$gen1.Range("G4").Value = "['001000' '110110']"
$gen1.Range("I4").Value = 0.5238095238095237
$gen1.Range("J4").Value = 48
$gen1.Range("K4").Value = "110000"
$gen1.Range("L4").Value = 0.27437641723356

$gen1.Range("A5").Value = -0.746031746031746
$gen1.Range("B5").Value = 8
$gen1.Range("C5").Value = "001000"
$gen1.Range("D5").Value = 0.5565633660871756
$gen1.Range("F5").Value = "001110"
$gen1.Range("G5").Value = "['001000' '110110']"
$gen1.Range("I5").Value = -0.5555555555555556
$gen1.Range("J5").Value = 14
$gen1.Range("K5").Value = "001110"
$gen1.Range("L5").Value = 0.308641975308642

$gen1.Range("A6").Value = -0.7142857142857143
$gen1.Range("B6").Value = 9
$gen1.Range("C6").Value = "001001"
$gen1.Range("D6").Value = 0.5102040816326531
$gen1.Range("F6").Value = "000110"
$gen1.Range("G6").Value = "['110110' '000100']"
$gen1.Range("I6").Value = -0.8095238095238095
$gen1.Range("J6").Value = 6
$gen1.Range("K6").Value = "000110"
$gen1.Range("L6").Value = 0.655328798185941

$gen1.Range("A7").Value = -0.7142857142857143
$gen1.Range("B7").Value = 9
$gen1.Range("C7").Value = "001001"
$gen1.Range("D7").Value = 0.5102040816326531
$gen1.Range("F7").Value = "110100"
$gen1.Range("G7").Value = "['110110' '000100']"
$gen1.Range("I7").Value = 0.6507936507936507
$gen1.Range("J7").Value = 52
$gen1.Range("K7").Value = "110100"
$gen1.Range("L7").Value = 0.4235323759133282

$gen1.Range("A8").Value = -0.7142857142857143
$gen1.Range("B8").Value = 9
$gen1.Range("C8").Value = "001001"
$gen1.Range("D8").Value = 0.5102040816326531
$gen1.Range("F8").Value = "001001"
$gen1.Range("G8").Value = "['001001' '001000']"
$gen1.Range("I8").Value = -0.7142857142857143
$gen1.Range("J8").Value = 9
$gen1.Range("K8").Value = "001001"
$gen1.Range("L8").Value = 0.5102040816326531

$gen1.Range("A9").Value = -0.7142857142857143
$gen1.Range("B9").Value = 9
$gen1.Range("C9").Value = "001001"
$gen1.Range("D9").Value = 0.5102040816326531
$gen1.Range("F9").Value = "001000"
$gen1.Range("G9").Value = "['001001' '001000']"
$gen1.Range("I9").Value = -0.746031746031746
$gen1.Range("J9").Value = 8
$gen1.Range("K9").Value = "001000"
$gen1.Range("L9").Value = 0.5565633660871756

$gen1.Range("F10").Value = "001110"
$gen1.Range("G10").Value = "['110110' '001000']"
$gen1.Range("I10").Value = -0.5555555555555556
$gen1.Range("J10").Value = 14
$gen1.Range("K10").Value = "001110"
$gen1.Range("L10").Value = 0.308641975308642

$gen1.Range("A11").Value = 0.6507936507936507
$gen1.Range("B11").Value = 52
$gen1.Range("C11").Value = "110100"
$gen1.Range("D11").Value = 0.4235323759133282
$gen1.Range("F11").Value = "110000"
$gen1.Range("G11").Value = "['110110' '001000']"
$gen1.Range("I11").Value = 0.5238095238095237
$gen1.Range("J11").Value = 48
$gen1.Range("K11").Value = "110000"
$gen1.Range("L11").Value = 0.27437641723356

# ---------------------------------------------------------------------------
# 4. Update gen_2 (sheet3) rows 2-11
# ---------------------------------------------------------------------------
$gen2 = $wb.Worksheets.Item("gen_2")

$gen2.Range("A2").Value = -0.7142857142857143
$gen2.Range("B2").Value = 9
$gen2.Range("C2").Value = "001001"
$gen2.Range("D2").Value = 0.5102040816326531
$gen2.Range("F2").Value = "001001"
$gen2.Range("G2").Value = "['001001' '001000']"
$gen2.Range("I2").Value = -0.7142857142857143
$gen2.Range("J2").Value = 9
$gen2.Range("K2").Value = "001001"
$gen2.Range("L2").Value = 0.5102040816326531

$gen2.Range("A3").Value = -0.873015873015873
$gen2.Range("B3").Value = 4
$gen2.Range("C3").Value = "000100"
$gen2.Range("D3").Value = 0.762156714537667
$gen2.Range("F3").Value = "001000"
$gen2.Range("G3").Value = "['001001' '001000']"
$gen2.Range("I3").Value = -0.746031746031746
$gen2.Range("J3").Value = 8
$gen2.Range("K3").Value = "001000"
$gen2.Range("L3").Value = 0.5565633660871756

$gen2.Range("A4").Value = -0.8095238095238095
$gen2.Range("B4").Value = 6
$gen2.Range("C4").Value = "000110"
$gen2.Range("D4").Value = 0.655328798185941
$gen2.Range("F4").Value = "001000"
$gen2.Range("G4").Value = "['001000' '001000']"
$gen2.Range("I4").Value = -0.746031746031746
$gen2.Range("J4").Value = 8
$gen2.Range("K4").Value = "001000"
$gen2.Range("L4").Value = 0.5565633660871756

$gen2.Range("A5").Value = -0.746031746031746
$gen2.Range("B5").Value = 8
$gen2.Range("C5").Value = "001000"
$gen2.Range("D5").Value = 0.5565633660871756
$gen2.Range("F5").Value = "001000"
$gen2.Range("G5").Value = "['001000' '001000']"
$gen2.Range("I5").Value = -0.746031746031746
$gen2.Range("J5").Value = 8
$gen2.Range("K5").Value = "001000"
$gen2.Range("L5").Value = 0.5565633660871756

$gen2.Range("A6").Value = -0.746031746031746
$gen2.Range("B6").Value = 8
$gen2.Range("C6").Value = "001000"
$gen2.Range("D6").Value = 0.5565633660871756
$gen2.Range("F6").Value = "000000"
$gen2.Range("G6").Value = "['001000' '000110']"
$gen2.Range("I6").Value = -1
$gen2.Range("J6").Value = 0
$gen2.Range("K6").Value = "000000"
$gen2.Range("L6").Value = 1

$gen2.Range("A7").Value = -0.746031746031746
$gen2.Range("B7").Value = 8
$gen2.Range("C7").Value = "001000"
$gen2.Range("D7").Value = 0.5565633660871756
$gen2.Range("F7").Value = "001110"
$gen2.Range("G7").Value = "['001000' '000110']"
$gen2.Range("I7").Value = -0.5555555555555556
$gen2.Range("J7").Value = 14
$gen2.Range("K7").Value = "001110"
$gen2.Range("L7").Value = 0.308641975308642

$gen2.Range("A8").Value = -0.746031746031746
$gen2.Range("B8").Value = 8
$gen2.Range("C8").Value = "001000"
$gen2.Range("D8").Value = 0.5565633660871756
$gen2.Range("F8").Value = "001000"
$gen2.Range("G8").Value = "['001000' '001000']"
$gen2.Range("I8").Value = -0.746031746031746
$gen2.Range("J8").Value = 8
$gen2.Range("K8").Value = "001000"
$gen2.Range("L8").Value = 0.5565633660871756

$gen2.Range("A9").Value = -0.746031746031746
$gen2.Range("B9").Value = 8
$gen2.Range("C9").Value = "001000"
$gen2.Range("D9").Value = 0.5565633660871756
$gen2.Range("F9").Value = "001000"
$gen2.Range("G9").Value = "['001000' '001000']"
$gen2.Range("I9").Value = -0.746031746031746
$gen2.Range("J9").Value = 8
$gen2.Range("K9").Value = "001000"
$gen2.Range("L9").Value = 0.5565633660871756

$gen2.Range("A10").Value = -0.7142857142857143
$gen2.Range("B10").Value = 9
$gen2.Range("C10").Value = "001001"
$gen2.Range("D10").Value = 0.5102040816326531
$gen2.Range("F10").Value = "000000"
$gen2.Range("G10").Value = "['001000' '000001']"
$gen2.Range("I10").Value = -1
$gen2.Range("J10").Value = 0
$gen2.Range("K10").Value = "000000"
$gen2.Range("L10").Value = 1

$gen2.Range("A11").Value = -0.7142857142857143
$gen2.Range("B11").Value = 9
$gen2.Range("C11").Value = "001001"
$gen2.Range("D11").Value = 0.5102040816326531
$gen2.Range("F11").Value = "001001"
$gen2.Range("G11").Value = "['001000' '000001']"
$gen2.Range("I11").Value = -0.7142857142857143
$gen2.Range("J11").Value = 9
$gen2.Range("K11").Value = "001001"
$gen2.Range("L11").Value = 0.5102040816326531

# ---------------------------------------------------------------------------
# 5. Populate the new gen_3 sheet (copy the gen_0 layout/headers, then fill data)
# ---------------------------------------------------------------------------
$gen0.Range("A1:L1").Copy()
$gen3.Range("A1").PasteSpecial()

$gen3.Range("A2").Value = -1
$gen3.Range("B2").Value = 0
$gen3.Range("C2").Value = "000000"
$gen3.Range("D2").Value = 1
$gen3.Range("F2").Value = "000000"
$gen3.Range("G2").Value = "['000000' '000110']"
$gen3.Range("I2").Value = -1
$gen3.Range("J2").Value = 0
$gen3.Range("K2").Value = "000000"
$gen3.Range("L2").Value = 1

$gen3.Range("A3").Value = -1
$gen3.Range("B3").Value = 0
$gen3.Range("C3").Value = "000000"
$gen3.Range("D3").Value = 1
$gen3.Range("F3").Value = "000110"
$gen3.Range("G3").Value = "['000000' '000110']"
$gen3.Range("I3").Value = -0.8095238095238095
$gen3.Range("J3").Value = 6
$gen3.Range("K3").Value = "000110"
$gen3.Range("L3").Value = 0.655328798185941

$gen3.Range("A4").Value = -0.9682539682539683
$gen3.Range("B4").Value = 1
$gen3.Range("C4").Value = "000001"
$gen3.Range("D4").Value = 0.9375157470395565
$gen3.Range("F4").Value = "000000"
$gen3.Range("G4").Value = "['001000' '000000']"
$gen3.Range("I4").Value = -1
$gen3.Range("J4").Value = 0
$gen3.Range("K4").Value = "000000"
$gen3.Range("L4").Value = 1

$gen3.Range("A5").Value = -0.873015873015873
$gen3.Range("B5").Value = 4
$gen3.Range("C5").Value = "000100"
$gen3.Range("D5").Value = 0.762156714537667
$gen3.Range("F5").Value = "001000"
$gen3.Range("G5").Value = "['001000' '000000']"
$gen3.Range("I5").Value = -0.746031746031746
$gen3.Range("J5").Value = 8
$gen3.Range("K5").Value = "001000"
$gen3.Range("L5").Value = 0.5565633660871756

$gen3.Range("A6").Value = -0.8095238095238095
$gen3.Range("B6").Value = 6
$gen3.Range("C6").Value = "000110"
$gen3.Range("D6").Value = 0.655328798185941
$gen3.Range("F6").Value = "000000"
$gen3.Range("G6").Value = "['001000' '000000']"
$gen3.Range("I6").Value = -1
$gen3.Range("J6").Value = 0
$gen3.Range("K6").Value = "000000"
$gen3.Range("L6").Value = 1

$gen3.Range("A7").Value = -0.746031746031746
$gen3.Range("B7").Value = 8
$gen3.Range("C7").Value = "001000"
$gen3.Range("D7").Value = 0.5565633660871756
$gen3.Range("F7").Value = "001000"
$gen3.Range("G7").Value = "['001000' '000000']"
$gen3.Range("I7").Value = -0.746031746031746
$gen3.Range("J7").Value = 8
$gen3.Range("K7").Value = "001000"
$gen3.Range("L7").Value = 0.5565633660871756

$gen3.Range("A8").Value = -0.746031746031746
$gen3.Range("B8").Value = 8
$gen3.Range("C8").Value = "001000"
$gen3.Range("D8").Value = 0.5565633660871756
$gen3.Range("F8").Value = "000000"
$gen3.Range("G8").Value = "['000000' '000001']"
$gen3.Range("I8").Value = -1
$gen3.Range("J8").Value = 0
$gen3.Range("K8").Value = "000000"
$gen3.Range("L8").Value = 1

$gen3.Range("A9").Value = -0.746031746031746
$gen3.Range("B9").Value = 8
$gen3.Range("C9").Value = "001000"
$gen3.Range("D9").Value = 0.5565633660871756
$gen3.Range("F9").Value = "000001"
$gen3.Range("G9").Value = "['000000' '000001']"
$gen3.Range("I9").Value = -0.9682539682539683
$gen3.Range("J9").Value = 1
$gen3.Range("K9").Value = "000001"
$gen3.Range("L9").Value = 0.9375157470395565

$gen3.Range("A10").Value = -0.746031746031746
$gen3.Range("B10").Value = 8
$gen3.Range("C10").Value = "001000"
$gen3.Range("D10").Value = 0.5565633660871756
$gen3.Range("F10").Value = "001000"
$gen3.Range("G10").Value = "['000000' '001000']"
$gen3.Range("I10").Value = -0.746031746031746
$gen3.Range("J10").Value = 8
$gen3.Range("K10").Value = "001000"
$gen3.Range("L10").Value = 0.5565633660871756

$gen3.Range("A11").Value = -0.746031746031746
$gen3.Range("B11").Value = 8
$gen3.Range("C11").Value = "001000"
$gen3.Range("D11").Value = 0.5565633660871756
$gen3.Range("F11").Value = "000000"
$gen3.Range("G11").Value = "['000000' '001000']"
$gen3.Range("I11").Value = -1
$gen3.Range("J11").Value = 0
$gen3.Range("K11").Value = "000000"
$gen3.Range("L11").Value = 1

# ---------------------------------------------------------------------------
# 6. Update "Función objetivo" sheet rows 8 & 11
# ---------------------------------------------------------------------------
$funcSheet.Range("A8").Value = -1
$funcSheet.Range("B8").Value = 0
$funcSheet.Range("C8").Value = "000000"
$funcSheet.Range("D8").Value = 1

$funcSheet.Range("A11").Value = -0.873015873015873
$funcSheet.Range("B11").Value = 4
$funcSheet.Range("C11").Value = "000100"
$funcSheet.Range("D11").Value = 0.762156714537667
